# Updates on 12th Oct
# Adds four new worksheets (AdminSettings_User, AdminSettings_Team,
# Team_PageDesign_AppLanding, Team_PageDesign_Onboarding) with their test
# data, matching the target commit.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write a whole header/data row from an array of values, skipping
# any $null entries (so we don't clobber cells that should stay blank).
# ---------------------------------------------------------------------------
function Set-RowValues {
    param($ws, [int]$row, [object[]]$values)
    for ($i = 0; $i -lt $values.Length; $i++) {
        $val = $values[$i]
        if ($null -ne $val) {
            $ws.Cells.Item($row, $i + 1).Value = $val
        }
    }
}

# ===========================================================================
# Sheet: AdminSettings_User
# ===========================================================================
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsUser = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsUser.Name = "AdminSettings_User"

Set-RowValues $wsUser 1 @("Description","TestType","FirstName","LastName","Email","Role","Success_or_Error?",$null,"Expected_Msg_FieldLevel","Expected_Msg_Header")
$wsUser.Range("A1:I1").Font.Bold = $true

Set-RowValues $wsUser 2 @("Select Role Empty","Negative","Ravi","B",$null,$null,"Error",$null,"Please select user's role")
$wsUser.Range("E2").Formula = '="Testinguser"&TEXT(NOW(),"ddmm")&"@gmail.com"'

Set-RowValues $wsUser 3 @("FirstName is Empty","Negative",$null,"B",$null,"User","Error",$null,"Please Enter First Name")
$wsUser.Range("E3").Formula = '="Testinguser"&TEXT(NOW(),"ddmm")&"@gmail.com"'

Set-RowValues $wsUser 4 @("Special chars in FirstName","Negative","r@v8","B",$null,"User","Error",$null,"Please enter letters only")
$wsUser.Range("E4").Formula = '="Testinguser"&TEXT(NOW(),"ddmm")&"@gmail.com"'

Set-RowValues $wsUser 5 @("Empty Email","Negative","Ravi","B",$null,"User","Error",$null,"Please Enter Email")

Set-RowValues $wsUser 6 @("Invalid Email","Negative","Ravi","B","rav","User","Error","Invalid Email")

Set-RowValues $wsUser 7 @("Successful Creation_User","Possitive","Ravi","B",$null,"User","Success",$null,"User created successfully. Email sent to user with credentials.")
$wsUser.Range("E7").Formula = '="Testinguser"&TEXT(NOW(),"ddmm")&"@gmail.com"'

Set-RowValues $wsUser 8 @("Successful Creation_SpAdmin","Possitive","Bolla","R",$null,"Specialty Admin","Success",$null,"User created successfully. Email sent to user with credentials.")
$wsUser.Range("E8").Formula = '="Testingspa"&TEXT(NOW(),"ddmm")&"@gmail.com"'

# ===========================================================================
# Sheet: AdminSettings_Team
# ===========================================================================
$wsTeam = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsUser)
$wsTeam.Name = "AdminSettings_Team"

Set-RowValues $wsTeam 1 @("Description","TestType","Speciality","TeamName",$null,"Expected_Msg_FieldLevel","Expected_Msg_Header","MemberEmail","SelectRole")
$wsTeam.Range("A1:H1").Font.Bold = $true

Set-RowValues $wsTeam 2 @("Create New Team","Positive","Testing")
$wsTeam.Range("D2").Formula = '="TestTeam"&TEXT(NOW()+4,"ddmm")'
$wsTeam.Range("F2").Formula = '="Successfully created team : "&D2'

Set-RowValues $wsTeam 3 @("Create New Team with existing data","Negative","Testing")
$wsTeam.Range("D3").Formula = '="TestTeam"&TEXT(NOW()+4,"ddmm")'
Set-RowValues $wsTeam 3 @($null,$null,$null,$null,$null,"A team already exists with the same name")

Set-RowValues $wsTeam 4 @("Add members to the team","Positive")
$wsTeam.Range("C4").Formula = '=C2'
$wsTeam.Range("D4").Formula = '=D2'
Set-RowValues $wsTeam 4 @($null,$null,$null,$null,$null,"Successfully added user to the team","tmsrafi.qa@gmail.com","Clinician")

$wsTeam.Range("G5").Formula = '=AdminSettings_User!E7'

# ===========================================================================
# Sheet: Team_PageDesign_AppLanding
# ===========================================================================
$wsAppLanding = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsTeam)
$wsAppLanding.Name = "Team_PageDesign_AppLanding"

Set-RowValues $wsAppLanding 1 @("Description","TestType","BrandColor","FontColor","Font","BrandLogo","AppTagline","AppDesc",$null,"Expected_Msg_Header")
$wsAppLanding.Range("A1:I1").Font.Bold = $true

Set-RowValues $wsAppLanding 2 @("Page Design in App Landing","Positive","#dddd34","Black","Open Sans",$null,"Video Consultation","Get medical advice and assistance without the need to visit the hospital.",$null,"Successfully updated team page design")
$wsAppLanding.Range("F2").Value = "D:\IMI_Automation\IMIAssist_Automtion\TestData\NewLogo.jpeg"
$wsAppLanding.Range("F2").Interior.Color = 65535

# ===========================================================================
# Sheet: Team_PageDesign_Onboarding
# ===========================================================================
$wsOnboarding = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsAppLanding)
$wsOnboarding.Name = "Team_PageDesign_Onboarding"

Set-RowValues $wsOnboarding 1 @("Description","TestType","ButtonColor","ButtonLabel","Expected_Msg_Header")
$wsOnboarding.Range("A1:E1").Font.Bold = $true

Set-RowValues $wsOnboarding 2 @("Page Design Onboarding","Positive","#dddd34","Proceed to Call.","Successfully updated team page design")

$wsOnboarding.Activate()
$wsOnboarding.Range("A2").Select()
